# Generate Report for Handback
# Update the Correspond Handoff/Handback Datetime values on the zh-cn and
# de-de sheets to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-30 10:43:43"
$wsZhCn.Range("H2").Value = "2016-03-30 10:44:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-30 10:43:55"
$wsDeDe.Range("H2").Value = "2016-03-30 10:45:01"
